$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New role names appended to the position table (column C), with
# departmentID (column B) = 1, continuing the id sequence in column A.
$newRoles = @(
    "Engineering Manager",
    "Program Manager",
    "DQA",
    "Subject Matter Expert",
    "Principal Technologist",
    "SQM",
    "Principal Consultant"
)

$startRow = 7
$startId = 6

for ($i = 0; $i -lt $newRoles.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $startId + $i
    $ws.Cells.Item($row, 2).Value = 1
    $ws.Cells.Item($row, 2).HorizontalAlignment = -4108
    $ws.Cells.Item($row, 3).Value = $newRoles[$i]
}

# Update the active selection to match the authored state.
$ws.Range("E12").Select()
